$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.947.76"
$ws.Range("E2").Value = "  -0.36%  "
Set-TextValue $ws.Range("D3") "1.823.38"
$ws.Range("E3").Value = "  -0.44%  "
Set-TextValue $ws.Range("D4") "0.9959"
$ws.Range("E4").Value = "  -0.33%  "
Set-TextValue $ws.Range("D5") "243.34"
$ws.Range("E5").Value = "  +0.74%  "
Set-TextValue $ws.Range("D6") "0.6297"
$ws.Range("E6").Value = "  +0.40%  "
Set-TextValue $ws.Range("D7") "0.9977"
$ws.Range("E7").Value = "  -0.28%  "
Set-TextValue $ws.Range("D8") "0.07480"
$ws.Range("E8").Value = "  -1.51%  "
Set-TextValue $ws.Range("D9") "0.2931"
$ws.Range("E9").Value = "  +0.42%  "
Set-TextValue $ws.Range("D10") "22.98"
$ws.Range("E10").Value = "  +0.56%  "
Set-TextValue $ws.Range("D11") "0.07684"
$ws.Range("E11").Value = "  +0.56%  "
Set-TextValue $ws.Range("D12") "1.829.88"
$ws.Range("E12").Value = "  -0.09%  "
Set-TextValue $ws.Range("D13") "4.981"
$ws.Range("E13").Value = "  +0.50%  "
Set-TextValue $ws.Range("D14") "0.6656"
$ws.Range("E14").Value = "  +0.11%  "
Set-TextValue $ws.Range("D15") "82.84"
$ws.Range("E15").Value = "  +0.54%  "
Set-TextValue $ws.Range("D16") "0.000009581"
$ws.Range("E16").Value = "  +1.83%  "
Set-TextValue $ws.Range("D17") "6.032"
$ws.Range("E17").Value = "  +0.85%  "
Set-TextValue $ws.Range("D18") "28.991.26"
$ws.Range("E18").Value = "  +0.08%  "
Set-TextValue $ws.Range("D19") "12.54"
$ws.Range("E19").Value = "  +1.75%  "
Set-TextValue $ws.Range("D20") "225.36"
$ws.Range("E20").Value = "  +0.19%  "
Set-TextValue $ws.Range("D21") "0.9962"
$ws.Range("E21").Value = "  -0.34%  "
Set-TextValue $ws.Range("D22") "7.118"
$ws.Range("E22").Value = "  -1.38%  "
Set-TextValue $ws.Range("D23") "0.9970"
$ws.Range("E23").Value = "  -0.39%  "
Set-TextValue $ws.Range("D24") "159.97"
$ws.Range("E24").Value = "  -0.71%  "
Set-TextValue $ws.Range("D25") "0.1410"
$ws.Range("E25").Value = "  +3.54%  "
Set-TextValue $ws.Range("D26") "8.479"
$ws.Range("E26").Value = "  +0.75%  "
Set-TextValue $ws.Range("D27") "17.88"
$ws.Range("E27").Value = "  +0.23%  "
Set-TextValue $ws.Range("D28") "1.497"
$ws.Range("E28").Value = "  +0.19%  "
Set-TextValue $ws.Range("D29") "4.120"
$ws.Range("E29").Value = "  +1.53%  "
Set-TextValue $ws.Range("D30") "4.040"
$ws.Range("E30").Value = "  +0.06%  "
Set-TextValue $ws.Range("D31") "0.05424"
$ws.Range("E31").Value = "  +4.27%  "
Set-TextValue $ws.Range("D32") "1.197"
$ws.Range("E32").Value = "  +0.09%  "
Set-TextValue $ws.Range("D33") "1.850"
$ws.Range("E33").Value = "  -0.02%  "
Set-TextValue $ws.Range("D34") "0.7403"
$ws.Range("E34").Value = "  +1.49%  "
Set-TextValue $ws.Range("D35") "1.132"
$ws.Range("E35").Value = "  -1.66%  "
Set-TextValue $ws.Range("D36") "2.634"
$ws.Range("E36").Value = "  +1.52%  "
Set-TextValue $ws.Range("D37") "1.235.84"
$ws.Range("E37").Value = "  -2.89%  "
Set-TextValue $ws.Range("D38") "2.745"
$ws.Range("E38").Value = "  -0.65%  "
Set-TextValue $ws.Range("D39") "0.01775"
$ws.Range("E39").Value = "  -0.61%  "
Set-TextValue $ws.Range("D40") "6.636"
$ws.Range("E40").Value = "  +2.05%  "
Set-TextValue $ws.Range("D41") "0.8979"
$ws.Range("E41").Value = "  +0.75%  "
Set-TextValue $ws.Range("D42") "0.9974"
$ws.Range("E42").Value = "  -0.29%  "
Set-TextValue $ws.Range("D43") "101.11"
$ws.Range("E43").Value = "  -0.46%  "
Set-TextValue $ws.Range("D44") "1.977.81"
$ws.Range("E44").Value = "  +0.17%  "
Set-TextValue $ws.Range("D45") "0.00000000123"
$ws.Range("E45").Value = "  +2.18%  "
Set-TextValue $ws.Range("D46") "64.95"
$ws.Range("E46").Value = "  +2.02%  "
Set-TextValue $ws.Range("D47") "0.5081"
$ws.Range("E47").Value = "  -0.51%  "
Set-TextValue $ws.Range("D48") "0.4034"
$ws.Range("E48").Value = "  +1.29%  "
Set-TextValue $ws.Range("D49") "8.934"
$ws.Range("E49").Value = "  +1.60%  "
Set-TextValue $ws.Range("D50") "0.07243"
$ws.Range("E50").Value = "  -1.16%  "

# Row 51: Cronos -> RenderToken
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D51") "1.650"
$ws.Range("E51").Value = "  +0.45%  "

